$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 header cells: labels were reworked so the "expiry" / "lot" labels
# become short static prefixes ("No.Lot:" / "Exp:") and the previously
# hard-coded values move into their own adjacent number cells.
$ws.Range("C3").Value = "utd"
$ws.Range("H3").Value = "No.Lot:"
$ws.Range("I3").Value = 5235325
$ws.Range("P3").Value = "Exp:"
$ws.Range("Q3").Value = 287652

# Restore the active selection to G12 on Sheet1.
$ws.Range("G12").Select()
